$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts them to numeric values
# (the source data models every cell in this sheet as text).
$textCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D14", "D15", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D37", "D38", "D39", "D41", "D43", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "41.982.69"
$ws.Range("E2").Value = "  -2.14%  "
$ws.Range("D3").Value = "2.289.96"
$ws.Range("E3").Value = "  -3.11%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "317.99"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").Value = "103.55"
$ws.Range("E6").Value = "  -3.44%  "
$ws.Range("E7").Value = "  -0.91%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "0.605"
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("D10").Value = "39.35"
$ws.Range("E10").Value = "  -5.55%  "
$ws.Range("D11").Value = "0.0904"
$ws.Range("E11").Value = "  -2.49%  "
$ws.Range("D12").Value = "8.25"
$ws.Range("E12").Value = "  -2.90%  "
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("D14").Value = "0.961"
$ws.Range("E14").Value = "  -4.46%  "
$ws.Range("D15").Value = "15.22"
$ws.Range("E15").Value = "  -5.29%  "
$ws.Range("D16").Value = "2.636.90"
$ws.Range("E16").Value = "  -3.01%  "
$ws.Range("D17").Value = "2.293.81"
$ws.Range("E17").Value = "  -2.25%  "
$ws.Range("D18").Value = "42.049.53"
$ws.Range("E18").Value = "  -1.90%  "
$ws.Range("D19").Value = "7.42"
$ws.Range("E19").Value = "  -2.24%  "
$ws.Range("D20").Value = "0.0000105"
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("B21").Value = "PancakeSwap"
$ws.Range("C21").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D21").Value = "3.63"
$ws.Range("E21").Value = "  -1.30%  "
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").Value = "73.30"
$ws.Range("E22").Value = "  -3.69%  "
$ws.Range("D23").Value = "279.53"
$ws.Range("E23").Value = "  +5.09%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "10.04"
$ws.Range("E24").Value = "  +6.61%  "
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").Value = "2.26"
$ws.Range("E25").Value = "  -2.67%  "
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("D27").Value = "10.79"
$ws.Range("E27").Value = "  -5.41%  "
$ws.Range("D28").Value = "2.39"
$ws.Range("E28").Value = "  +6.05%  "
$ws.Range("D29").Value = "22.91"
$ws.Range("E29").Value = "  -2.75%  "
$ws.Range("D30").Value = "36.13"
$ws.Range("E30").Value = "  -1.94%  "
$ws.Range("D31").Value = "163.01"
$ws.Range("E31").Value = "  -4.24%  "
$ws.Range("D32").Value = "0.0870"
$ws.Range("E32").Value = "  -3.06%  "
$ws.Range("E33").Value = "  -2.39%  "
$ws.Range("D34").Value = "5.82"
$ws.Range("E34").Value = "  -3.15%  "
$ws.Range("E35").Value = "  +4.15%  "
$ws.Range("E36").Value = "  -4.77%  "
$ws.Range("D37").Value = "4.51"
$ws.Range("E37").Value = "  -4.43%  "
$ws.Range("D38").Value = "2.91"
$ws.Range("E38").Value = "  +7.56%  "
$ws.Range("D39").Value = "0.0349"
$ws.Range("E39").Value = "  -3.80%  "
$ws.Range("E40").Value = "  -3.48%  "
$ws.Range("D41").Value = "99.21"
$ws.Range("E41").Value = "  -3.89%  "
$ws.Range("E42").Value = "  -4.73%  "
$ws.Range("D43").Value = "69.41"
$ws.Range("E43").Value = "  -2.48%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("E45").Value = "  -6.25%  "
$ws.Range("D46").Value = "112.86"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D47").Value = "11.86"
$ws.Range("E47").Value = "  -3.29%  "
$ws.Range("D48").Value = "76.74"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").Value = "8.97"
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("D50").Value = "5.28"
$ws.Range("E50").Value = "  -4.87%  "
$ws.Range("D51").Value = "1.578.58"
$ws.Range("E51").Value = "  +0.05%  "
